# Applies the "Updated cryptos list" data refresh to sheet1 (crypto price/volume table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "79.788.36"
$ws.Range("E2").Value = "  +4.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.202.92"
$ws.Range("E3").Value = "  +4.95%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.20"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "633.88"
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.238"
$ws.Range("E8").Value = "  +13.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  +5.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.199.45"
$ws.Range("E10").Value = "  +4.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +32.39%  "

$ws.Range("E12").Value = "  +3.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.50"
$ws.Range("E13").Value = "  +6.91%  "

$ws.Range("E14").Value = "  +19.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.793.93"
$ws.Range("E15").Value = "  +4.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.95"
$ws.Range("E16").Value = "  +8.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.642.03"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.202.72"
$ws.Range("E18").Value = "  +5.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.52"
$ws.Range("E19").Value = "  +6.94%  "

$ws.Range("E20").Value = "  +29.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.18"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.09"
$ws.Range("E22").Value = "  +14.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.11"
$ws.Range("E23").Value = "  +17.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +13.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.369.54"
$ws.Range("E25").Value = "  +5.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "77.32"
$ws.Range("E26").Value = "  +5.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.74"
$ws.Range("E27").Value = "  +6.95%  "

$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  +6.89%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.03"
$ws.Range("E31").Value = "  +8.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").Value = "  +5.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "527.49"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.99"
$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.142"
$ws.Range("E35").Value = "  +25.45%  "

$ws.Range("E36").Value = "  +10.21%  "

$ws.Range("E37").Value = "  +11.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.407"
$ws.Range("E39").Value = "  +5.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.34"
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.05"
$ws.Range("E41").Value = "  +0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "192.84"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  +6.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.821"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("E46").Value = "  +7.10%  "

$ws.Range("E47").Value = "  +3.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.32"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.96"
$ws.Range("E49").Value = "  +15.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.639"
$ws.Range("E50").Value = "  +4.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.51"
$ws.Range("E51").Value = "  +1.26%  "

# Restore default (unstyled) formatting on the Price cells we touched, so only the
# underlying text values change and no new cell formatting is introduced.
$ws.Range("D2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
